# Natmi following Dr Hou advice
# Rebuild the sending/target cluster matrix (ECs, FAPs, M2, sCs) x (ECs, FAPs, M2, sCs)
# for the Fn1-Itga2 ligand-receptor pair, replacing the previous partial table
# (rows 2-9) with the complete 16-row table (rows 2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value, K Receptor-expressing cells,
# L Receptor detection rate, M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value, P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$rows = @(
        @("ECs", "Fn1", "Itga2", "ECs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 2, 0.6666666666666666, 2.672731, 8.018193, 0.5408378022089502, 0.5408378022089503, 108.9228523717313, 980.305671345582, 0.01150111397927924, 0.01150111397927925),
        @("ECs", "Fn1", "Itga2", "FAPs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 3, 1, 1.627877666666667, 4.883633, 0.3294075533620981, 0.3294075533620982, 66.34153559246022, 597.073820332142, 0.007004972288141408, 0.00700497228814141),
        @("ECs", "Fn1", "Itga2", "M2", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 1, 0.3333333333333333, 0.01852966666666667, 0.055589, 0.003749552123152102, 0.003749552123152104, 0.7551467569428889, 6.796320812486, 0.00007973559940427397, 0.000079735599404274),
        @("ECs", "Fn1", "Itga2", "sCs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 3, 1, 0.6226963333333334, 1.868089, 0.1260050923057995, 0.1260050923057995, 25.37698735416511, 228.392886187486, 0.002679544445043638, 0.002679544445043639),
        @("FAPs", "Fn1", "Itga2", "ECs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 2, 0.6666666666666666, 2.672731, 8.018193, 0.5408378022089502, 0.5408378022089503, 4515.015897896506, 40635.14308106856, 0.4767384559738383, 0.4767384559738385),
        @("FAPs", "Fn1", "Itga2", "FAPs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 3, 1, 1.627877666666667, 4.883633, 0.3294075533620981, 0.3294075533620982, 2749.956334861484, 24749.60701375336, 0.2903666269897574, 0.2903666269897575),
        @("FAPs", "Fn1", "Itga2", "M2", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 1, 0.3333333333333333, 0.01852966666666667, 0.055589, 0.003749552123152102, 0.003749552123152104, 31.30196775609777, 281.71770980488, 0.003305160405733523, 0.003305160405733525),
        @("FAPs", "Fn1", "Itga2", "sCs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 3, 1, 0.6226963333333334, 1.868089, 0.1260050923057995, 0.1260050923057995, 1051.914257200542, 9467.22831480488, 0.1110711435209544, 0.1110711435209544),
        @("M2", "Fn1", "Itga2", "ECs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 2, 0.6666666666666666, 2.672731, 8.018193, 0.5408378022089502, 0.5408378022089503, 269.8533607437723, 2428.680246693951, 0.02849369248074487, 0.02849369248074488),
        @("M2", "Fn1", "Itga2", "FAPs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 3, 1, 1.627877666666667, 4.883633, 0.3294075533620981, 0.3294075533620982, 164.3593235644479, 1479.233912080031, 0.01735462552358337, 0.01735462552358337),
        @("M2", "Fn1", "Itga2", "M2", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 1, 0.3333333333333333, 0.01852966666666667, 0.055589, 0.003749552123152102, 0.003749552123152104, 1.870855250102555, 16.837697250923, 0.0001975427470144615, 0.0001975427470144616),
        @("M2", "Fn1", "Itga2", "sCs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 3, 1, 0.6226963333333334, 1.868089, 0.1260050923057995, 0.1260050923057995, 62.87078582649145, 565.8370724384231, 0.006638497413651954, 0.006638497413651956),
        @("sCs", "Fn1", "Itga2", "ECs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 2, 0.6666666666666666, 2.672731, 8.018193, 0.5408378022089502, 0.5408378022089503, 228.2852975929674, 2054.567678336706, 0.02410453977508766, 0.02410453977508767),
        @("sCs", "Fn1", "Itga2", "FAPs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 3, 1, 1.627877666666667, 4.883633, 0.3294075533620981, 0.3294075533620982, 139.0415038325762, 1251.373534493186, 0.01468132856061592, 0.01468132856061593),
        @("sCs", "Fn1", "Itga2", "M2", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 1, 0.3333333333333333, 0.01852966666666667, 0.055589, 0.003749552123152102, 0.003749552123152104, 1.582669737170889, 14.244027634538, 0.0001671133709998434, 0.0001671133709998435),
        @("sCs", "Fn1", "Itga2", "sCs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 3, 1, 0.6226963333333334, 1.868089, 0.1260050923057995, 0.1260050923057995, 53.18620458439312, 478.6758412595381, 0.005615906926149537, 0.005615906926149538)
)

$rowIndex = 2
foreach ($row in $rows) {
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
